$wb = $excel.ActiveWorkbook

# ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 869.9286
$ws.Range("I107").Value = 869.9286
$ws.Range("K107").Value = 869.9286
$ws.Range("M107").Value = 1050.0714
$ws.Range("H112").Value = 1716.3334
$ws.Range("I112").Value = 950
$ws.Range("J112").Value = 1812.125
$ws.Range("K112").Value = 2850
$ws.Range("L112").Value = 5436.375
$ws.Range("M112").Value = -1742
$ws.Range("N112").Value = -7652.375
$ws.Range("H129").Value = 898.4299999999999
$ws.Range("I129").Value = 607.5
$ws.Range("J129").Value = 917
$ws.Range("K129").Value = 1822.5
$ws.Range("L129").Value = 2751
$ws.Range("M129").Value = 3177.5
$ws.Range("N129").Value = -12751
$ws.Range("H132").Value = 3411.6086
$ws.Range("I132").Value = 3107.139
$ws.Range("J132").Value = 4507.7
$ws.Range("K132").Value = 9321.417000000001
$ws.Range("L132").Value = 13523.1
$ws.Range("M132").Value = -6791.417000000001
$ws.Range("N132").Value = -18583.1
$ws.Range("H137").Value = 29538.158
$ws.Range("I137").Value = 1384
$ws.Range("J137").Value = 68250.125
$ws.Range("K137").Value = 4152
$ws.Range("L137").Value = 204750.375
$ws.Range("M137").Value = -1602
$ws.Range("N137").Value = -209850.375
$ws.Range("H138").Value = 1548.625
$ws.Range("I138").Value = 894.5925999999999
$ws.Range("K138").Value = 2683.7778
$ws.Range("M138").Value = 2456.2222
$ws.Range("H139").Value = 14955.223
$ws.Range("J139").Value = 14955.223
$ws.Range("L139").Value = 14955.223
$ws.Range("N139").Value = -25235.223
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0

# ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1841.0605
$ws.Range("I61").Value = 1179.68
$ws.Range("J61").Value = 3907.875
$ws.Range("K61").Value = 1179.68
$ws.Range("L61").Value = 3907.875
$ws.Range("M61").Value = -967.6800000000001
$ws.Range("N61").Value = -4331.875
$ws.Range("H74").Value = 3312.775
$ws.Range("I74").Value = 3472.1892
$ws.Range("K74").Value = 3472.1892
$ws.Range("M74").Value = -2598.1892
$ws.Range("H77").Value = 3312.775
$ws.Range("I77").Value = 3472.1892
$ws.Range("K77").Value = 17360.946
$ws.Range("M77").Value = -12992.946
$ws.Range("H110").Value = 2858.6875
$ws.Range("I110").Value = 2185.3635
$ws.Range("J110").Value = 4340
$ws.Range("K110").Value = 2185.3635
$ws.Range("L110").Value = 4340
$ws.Range("M110").Value = -140.3634999999999
$ws.Range("N110").Value = -8430
$ws.Range("H136").Value = 1841.0605
$ws.Range("I136").Value = 1179.68
$ws.Range("J136").Value = 3907.875
$ws.Range("K136").Value = 3539.04
$ws.Range("L136").Value = 11723.625
$ws.Range("M136").Value = -989.04
$ws.Range("N136").Value = -16823.625

# BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 374.76923
$ws.Range("I94").Value = 316.2
$ws.Range("J94").Value = 570
$ws.Range("K94").Value = 316.2
$ws.Range("L94").Value = 570
$ws.Range("M94").Value = 134.8
$ws.Range("N94").Value = -1472
$ws.Range("H107").Value = 2308.32
$ws.Range("I107").Value = 2519.1
$ws.Range("J107").Value = 1465.2
$ws.Range("K107").Value = 2519.1
$ws.Range("L107").Value = 1465.2
$ws.Range("M107").Value = -599.0999999999999
$ws.Range("N107").Value = -5305.2
$ws.Range("H134").Value = 2231.5405
$ws.Range("I134").Value = 1966.2903
$ws.Range("K134").Value = 5898.8709
$ws.Range("M134").Value = -3363.8709

# CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4377.5938
$ws.Range("I31").Value = 3031.9
$ws.Range("J31").Value = 4989.273
$ws.Range("K31").Value = 3031.9
$ws.Range("L31").Value = 4989.273
$ws.Range("M31").Value = -2736.9
$ws.Range("N31").Value = -5579.273
$ws.Range("H34").Value = 4377.5938
$ws.Range("I34").Value = 3031.9
$ws.Range("J34").Value = 4989.273
$ws.Range("K34").Value = 3031.9
$ws.Range("L34").Value = 4989.273
$ws.Range("M34").Value = -2829.9
$ws.Range("N34").Value = -5393.273
$ws.Range("H132").Value = 1846.4117
$ws.Range("I132").Value = 1991.6061
$ws.Range("J132").Value = 1580.2222
$ws.Range("K132").Value = 5974.8183
$ws.Range("L132").Value = 4740.6666
$ws.Range("M132").Value = -3444.8183
$ws.Range("N132").Value = -9800.6666

# CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 715.12
$ws.Range("I5").Value = 651.5789
$ws.Range("J5").Value = 916.3333
$ws.Range("K5").Value = 1954.7367
$ws.Range("L5").Value = 2748.9999
$ws.Range("M5").Value = -1842.7367
$ws.Range("N5").Value = -2972.9999
$ws.Range("H92").Value = 557.1429000000001
$ws.Range("I92").Value = 1100
$ws.Range("J92").Value = 340
$ws.Range("K92").Value = 3300
$ws.Range("L92").Value = 1020
$ws.Range("M92").Value = -2052
$ws.Range("N92").Value = -3516
$ws.Range("H131").Value = 7649.067
$ws.Range("I131").Value = 564.1429000000001
$ws.Range("J131").Value = 13848.375
$ws.Range("K131").Value = 1692.4287
$ws.Range("L131").Value = 41545.125
$ws.Range("M131").Value = 3347.5713
$ws.Range("N131").Value = -51625.125
$ws.Range("H135").Value = 715.12
$ws.Range("I135").Value = 651.5789
$ws.Range("J135").Value = 916.3333
$ws.Range("K135").Value = 5864.2101
$ws.Range("L135").Value = 8246.9997
$ws.Range("M135").Value = -3329.2101
$ws.Range("N135").Value = -13316.9997

# LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1970
$ws.Range("I136").Value = 1476.0377
$ws.Range("J136").Value = 4151.6665
$ws.Range("K136").Value = 4428.1131
$ws.Range("L136").Value = 12454.9995
$ws.Range("M136").Value = -1878.1131
$ws.Range("N136").Value = -17554.9995

# WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1333.172
$ws.Range("I132").Value = 1086.7183
$ws.Range("J132").Value = 2128.5454
$ws.Range("K132").Value = 3260.1549
$ws.Range("L132").Value = 6385.6362
$ws.Range("M132").Value = -730.1549
$ws.Range("N132").Value = -11445.6362
$ws.Range("H136").Value = 2705.3872
$ws.Range("I136").Value = 2921.5908
$ws.Range("J136").Value = 2176.889
$ws.Range("K136").Value = 8764.7724
$ws.Range("L136").Value = 6530.667
$ws.Range("M136").Value = -6214.7724
$ws.Range("N136").Value = -11630.667

# ALC row 140: clear N140 entirely (cell removed from the sheet data)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N140").ClearContents()
